# Rename the "issuanceDate" placeholder to "dateRequested" in the two
# places where it appears as its own run (split across "{", "issuanceDate",
# "}" runs: the header date stamp and the "Date:" line), while leaving the
# "Issued on: {issuanceDate}" occurrence (a single run) untouched.
#
# A plain Find/Replace merges the run that is rewritten with its
# neighbouring runs that share identical formatting, even though those
# neighbours ("{" and "}") were not part of the match. Performing the
# replace with TrackRevisions enabled and then Accept()-ing just the
# resulting revisions keeps every run boundary intact, matching how the
# original runs ("{", "issuanceDate"->"dateRequested", "}") stay separate.

$d = $word.ActiveDocument
$originalTrackRevisions = $d.TrackRevisions
$d.TrackRevisions = $true

$changed = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*{issuanceDate}*" -and $t -notlike "*Issued on:*") {
        $p.Range.Find.Execute("issuanceDate", $true, $false, $false, $false, $false, $true, 1, $false, "dateRequested", 2)
        $changed = $changed + 1
    }
}
Write-Output "Paragraphs updated: $changed"

$d.TrackRevisions = $originalTrackRevisions

for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions.Item($i).Accept()
}
